$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new rows (5 and 6), inheriting formatting from row 4 above ---
$ws.Range("A5").EntireRow.Insert()
$ws.Range("A6").EntireRow.Insert()

# --- Row 2: user_number changes from 201500001 to 201500005 (now a duplicate of new row 4) ---
$ws.Range("A2").Value2 = 201500005

# --- Row 3: user_number changes to 201500006; password + section get a fresh value ---
$ws.Range("A3").Value2 = 201500006
$ws.Range("E3").Value2 = "be9f57a7bbea5f7489e601db0cecffcfdd91e508"
$ws.Range("H3").Value2 = "X42"

# --- Row 4: becomes a duplicate of the (updated) row 2 data ---
$ws.Range("A4").Value2 = 201500005
$ws.Range("B4").Value2 = "Juan"
$ws.Range("C4").Value2 = "Cruz"
$ws.Range("D4").Value2 = "Dela"
$ws.Range("E4").Value2 = "406039a9fe75eb67ada58f6a06b9a72410fb86e8"
$ws.Range("F4").Value2 = "jdc@fit.edu.ph"
$ws.Range("G4").Value2 = "BSITWMA"
$ws.Range("H4").Value2 = "W41"
$ws.Range("I4").Value2 = "4th"

# --- Row 5: duplicate of the original row 3 data (pre-update hash/section), new number ---
$ws.Range("A5").Value2 = 201500006
$ws.Range("B5").Value2 = "Two"
$ws.Range("C5").Value2 = "Cruz"
$ws.Range("D5").Value2 = "Dela"
$ws.Range("E5").Value2 = "be9f57a7bbea5f7489e601db0cecffcfdd91e507"
$ws.Range("F5").Value2 = "tdc@fit.edu.ph"
$ws.Range("G5").Value2 = "BSITAGD"
$ws.Range("H5").Value2 = "X41"
$ws.Range("I5").Value2 = "Terminal"

# --- Row 6: another duplicate of the original row 3 data, new unique number ---
$ws.Range("A6").Value2 = 201500007
$ws.Range("B6").Value2 = "Two"
$ws.Range("C6").Value2 = "Cruz"
$ws.Range("D6").Value2 = "Dela"
$ws.Range("E6").Value2 = "be9f57a7bbea5f7489e601db0cecffcfdd91e507"
$ws.Range("F6").Value2 = "tdc@fit.edu.ph"
$ws.Range("G6").Value2 = "BSITAGD"
$ws.Range("H6").Value2 = "X41"
$ws.Range("I6").Value2 = "Terminal"

# --- Rebuild the e-mail hyperlinks in the new order ---
$ws.Range("F2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F6"), "mailto:tdc@fit.edu.ph") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "mailto:jdc@fit.edu.ph") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "mailto:jdc@fit.edu.ph") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), "mailto:tdc@fit.edu.ph") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "mailto:tdc@fit.edu.ph") | Out-Null

# --- UI state: selection moved to C12 ---
$ws.Range("C12").Select()
